# Apply weekly price update to the Pomelo sheet.
# The underlying data rows (2-8) get their Fecha/Volumen/Precio columns
# shuffled around (a re-sort by date effectively), while the other columns
# (market, product, variety, unit, origin, kg/unit) stay identical since
# they are the same on every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg)
$rows = @(
    @{ Row = 2;  D = 44181; M = 65; N = 9000;  O = 10000; P = 9462;  S = 676 },
    @{ Row = 3;  D = 44253; M = 90; N = 12000; O = 13000; P = 12667; S = 905 },
    @{ Row = 4;  D = 44172; M = 90; N = 8500;  O = 9000;  P = 8806;  S = 629 },
    @{ Row = 5;  D = 44232; M = 60; N = 11000; O = 12000; P = 11583; S = 827 },
    @{ Row = 6;  D = 44229; M = 55; N = 11000; O = 12000; P = 11364; S = 812 },
    @{ Row = 7;  D = 44216; M = 55; N = 11000; O = 12000; P = 11545; S = 825 },
    @{ Row = 8;  D = 44210; M = 70; N = 10000; O = 11000; P = 10357; S = 740 }
)

foreach ($r in $rows) {
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("M" + $r.Row).Value = $r.M
    $ws.Range("N" + $r.Row).Value = $r.N
    $ws.Range("O" + $r.Row).Value = $r.O
    $ws.Range("P" + $r.Row).Value = $r.P
    $ws.Range("S" + $r.Row).Value = $r.S
}
